$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LIST")

# A3: rename RO.FOU.001.CRE -> RO.FOU.001.CRE.01 and strip its special style
$ws.Range("A3").Value = "RO.FOU.001.CRE.01"
$ws.Range("A3").Style = "Normal"

# A4: used to hold AD.SEC.014.FON.01 - clear the value, keep the blank "text" style
$ws.Range("A4").ClearContents()
$ws.Range("A4").NumberFormat = "@"

# A5: new blank cell carrying the same blank "text" style as its neighbours
$ws.Range("A5").NumberFormat = "@"

# A8, A10, A12, A18: fully cleared (contents + formatting) so the cell disappears
$ws.Range("A8").Clear()
$ws.Range("A10").Clear()
$ws.Range("A12").Clear()
$ws.Range("A18").Clear()

$ws.Range("A8:A9").Select()
$ws.Application.ActiveCell = $ws.Range("A9")

$ws2 = $wb.Worksheets.Item("Feuil1")
$ws2.Range("F2").Select()
